# Insert 3 new price-report rows before the existing row 891.
# Excel's row insert shifts rows 891-976 down to 894-979 and carries the
# row-890 formatting (incl. the date-format style on column D) onto the
# freshly inserted rows, matching the target file's dimension A1:T979.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("891:893").Insert()

# New row 891: Especial quality, Peru origin, $/bandeja 10 kilos unit
$ws.Range("A891").Value = 7
$ws.Range("B891").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C891").Value = "Ñuble"
$ws.Range("D891").Value = 45106
$ws.Range("E891").Value = 16
$ws.Range("F891").Value = "Fruta"
$ws.Range("G891").Value = 100106
$ws.Range("H891").Value = "Oleaginosos"
$ws.Range("I891").Value = 100106002
$ws.Range("J891").Value = "Palta"
$ws.Range("K891").Value = "Hass"
$ws.Range("L891").Value = "Especial"
$ws.Range("M891").Value = 120
$ws.Range("N891").Value = 30000
$ws.Range("O891").Value = 30000
$ws.Range("P891").Value = 30000
$ws.Range("Q891").Value = "$/bandeja 10 kilos"
$ws.Range("R891").Value = "Perú"
$ws.Range("S891").Value = 3000
$ws.Range("T891").Value = 10

# New row 892: Primera quality, Peru origin, $/bandeja 10 kilos unit
$ws.Range("A892").Value = 7
$ws.Range("B892").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C892").Value = "Ñuble"
$ws.Range("D892").Value = 45106
$ws.Range("E892").Value = 16
$ws.Range("F892").Value = "Fruta"
$ws.Range("G892").Value = 100106
$ws.Range("H892").Value = "Oleaginosos"
$ws.Range("I892").Value = 100106002
$ws.Range("J892").Value = "Palta"
$ws.Range("K892").Value = "Hass"
$ws.Range("L892").Value = "Primera"
$ws.Range("M892").Value = 100
$ws.Range("N892").Value = 25000
$ws.Range("O892").Value = 25000
$ws.Range("P892").Value = 25000
$ws.Range("Q892").Value = "$/bandeja 10 kilos"
$ws.Range("R892").Value = "Perú"
$ws.Range("S892").Value = 2500
$ws.Range("T892").Value = 10

# New row 893: Segunda quality, Peru origin, $/bandeja 10 kilos unit
$ws.Range("A893").Value = 7
$ws.Range("B893").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C893").Value = "Ñuble"
$ws.Range("D893").Value = 45106
$ws.Range("E893").Value = 16
$ws.Range("F893").Value = "Fruta"
$ws.Range("G893").Value = 100106
$ws.Range("H893").Value = "Oleaginosos"
$ws.Range("I893").Value = 100106002
$ws.Range("J893").Value = "Palta"
$ws.Range("K893").Value = "Hass"
$ws.Range("L893").Value = "Segunda"
$ws.Range("M893").Value = 80
$ws.Range("N893").Value = 23000
$ws.Range("O893").Value = 23000
$ws.Range("P893").Value = 23000
$ws.Range("Q893").Value = "$/bandeja 10 kilos"
$ws.Range("R893").Value = "Perú"
$ws.Range("S893").Value = 2300
$ws.Range("T893").Value = 10
